$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '61.467.10'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -2.19%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.368.90'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -3.22%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '403.70'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -2.98%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '132.19'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +5.03%  '
$ws.Range('E7').Value = '  -1.55%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.665'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -2.94%  '
$ws.Range('E10').Value = '  -7.60%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '41.94'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.78%  '
$ws.Range('E12').Value = '  -1.84%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '3.889.30'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -3.44%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '19.75'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.44%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '8.35'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -3.19%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.374.89'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -2.63%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '61.441.56'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -1.98%  '
$ws.Range('E18').Value = '  -2.92%  '
$ws.Range('E19').Value = '  -0.09%  '
$ws.Range('E20').Value = '  -8.43%  '
$ws.Range('E21').Value = '  -5.34%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '84.29'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +2.56%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '312.45'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -1.13%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '12.66'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -3.44%  '
$ws.Range('E25').Value = '  -3.01%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '4.77'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +10.40%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '29.35'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -5.42%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '8.19'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +1.51%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.73'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -1.61%  '
$ws.Range('E30').Value = '  +2.48%  '
$ws.Range('E31').Value = '  -3.34%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.115'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -2.82%  '
$ws.Range('E33').Value = '  -0.27%  '
$ws.Range('E34').Value = '  -2.90%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '41.51'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -2.62%  '
$ws.Range('E36').Value = '  -3.86%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '51.65'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -1.43%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.999'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +0.30%  '
$ws.Range('E39').Value = '  -3.53%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.92'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -4.14%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '138.47'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +1.87%  '
$ws.Range('E42').Value = '  -2.27%  '
$ws.Range('E43').Value = '  -1.49%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.291'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +1.07%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.96'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +0.65%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '16.60'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -2.39%  '
$ws.Range('E47').Value = '  -1.10%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '21.28'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -3.54%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.113.76'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -4.28%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.31'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -6.14%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.87'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.57%  '
